$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference an unstyled cell so we can reset style on rewritten D cells
# after forcing a text NumberFormat (avoids Excel auto-converting
# numeric-looking strings like '552.11' into floating point numbers).
$plainStyle = $ws.Range("D7").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.991.51"
$ws.Range("D2").Style = $plainStyle
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.415.88"
$ws.Range("D3").Style = $plainStyle
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "552.11"
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = "  -0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.92"
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = "  -1.89%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.595"
$ws.Range("D8").Style = $plainStyle
$ws.Range("E8").Value = "  +3.73%  "
$ws.Range("E9").Value = "  -2.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.66"
$ws.Range("D10").Style = $plainStyle
$ws.Range("E10").Value = "  -2.55%  "
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.352"
$ws.Range("D12").Style = $plainStyle
$ws.Range("E12").Value = "  -2.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.18"
$ws.Range("D13").Style = $plainStyle
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.845.98"
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = "  -1.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.926.43"
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("E16").Value = "  -2.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.415.58"
$ws.Range("D17").Style = $plainStyle
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.27"
$ws.Range("D18").Style = $plainStyle
$ws.Range("E18").Value = "  -2.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.41"
$ws.Range("D19").Style = $plainStyle
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "327.75"
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = "  -2.53%  "
$ws.Range("E21").Value = "  -3.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.84"
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = "  +1.36%  "
$ws.Range("E24").Value = "  +3.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.55"
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = "  -0.82%  "
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0774"
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = "  -3.21%  "
$ws.Range("E29").Value = "  -2.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.88"
$ws.Range("D30").Style = $plainStyle
$ws.Range("E30").Value = "  -1.22%  "
$ws.Range("E31").Value = "  -4.53%  "
$ws.Range("E32").Value = "  -1.74%  "
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  -1.47%  "
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("E37").Value = "  -2.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.61"
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = "  -2.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "323.28"
$ws.Range("D39").Style = $plainStyle
$ws.Range("E39").Value = "  +1.97%  "
$ws.Range("E40").Value = "  -3.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.65"
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = "  -2.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "139.92"
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = "  -3.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0970"
$ws.Range("D43").Style = $plainStyle
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.53"
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("E45").Value = "  -2.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.576"
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0223"
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = "  -1.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.384"
$ws.Range("D48").Style = $plainStyle
$ws.Range("E48").Value = "  -4.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.04"
$ws.Range("D49").Style = $plainStyle
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("E50").Value = "  -5.59%  "
$ws.Range("E51").Value = "  -0.97%  "
